# Applies the TERM I -> TERM II title change, plus the marks-sheet
# corrections to the subject table (GEOGRAPHY/AGRICULTURE,
# ENTREPRENEURSHIP/GEOGRAPHY, SUBMATH, GENERAL PAPER blocks).

$d = $word.ActiveDocument

# --- Title: TERM I -> TERM II -------------------------------------------
$d.Content.Find.Execute("TERM I", $true, $false, $false, $false, $false, $true, 0, $false, "TERM II", 1) | Out-Null

# --- Helper to get the Nth cell (1-based) of a table row as a Range -----
function Get-RowCell($row, [int]$index) {
    $i = 0
    foreach ($cell in $row.Cells) {
        $i = $i + 1
        if ($i -eq $index) {
            return $cell
        }
    }
    return $null
}

# Replace the (unique, whole-cell) text of a cell, leaving other runs
# inside the cell untouched, by searching for the exact old text within
# the cell's own Range.
#
# NOTE: a table cell's native .Range object must be used for plain text
# assignment (.Text = ...), but must *not* be used as the object that
# Find.Execute() is invoked on directly -- in that case the search
# scope is not respected and the first match anywhere in the document
# can be hit instead of the match inside the cell. Re-wrapping the same
# [start,end) offsets via $d.Range(...) keeps Find properly scoped to
# the cell.
function Set-CellText($row, [int]$index, [string]$oldText, [string]$newText) {
    $cell = Get-RowCell $row $index
    if ($oldText -eq "") {
        $cell.Range.Text = $newText
    } else {
        $cr = $cell.Range
        $findRng = $d.Range($cr.Start, $cr.End)
        $findRng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 0, $false, $newText, 1) | Out-Null
    }
}

$t = $d.Tables.Item(1)

# --- Row 4: GEOGRAPHY / Paper 1 / 56 / P7 / (subject grade O) -----------
$row4 = $t.Rows.Item(4)
Set-CellText $row4 1 "GEOGRAPHY" "AGRICULTURE"
Set-CellText $row4 2 "Paper 1" ""
Set-CellText $row4 3 "56" ""
Set-CellText $row4 4 "P7" ""

# --- Row 5: Paper 2 / 74 / C4 --------------------------------------------
$row5 = $t.Rows.Item(5)
Set-CellText $row5 3 "74" "64.0"
Set-CellText $row5 4 "C4" "C6"

# --- Row 6: Paper 3 / 27 / F9 --------------------------------------------
$row6 = $t.Rows.Item(6)
Set-CellText $row6 2 "Paper 3" ""
Set-CellText $row6 3 "27" ""
Set-CellText $row6 4 "F9" ""

# --- Row 7: ENTREPRENEURSHIP / Paper 1 / 42 / P8 / (subject grade P8) ---
$row7 = $t.Rows.Item(7)
Set-CellText $row7 1 "ENTREPRENEURSHIP" "GEOGRAPHY"
Set-CellText $row7 2 "Paper 1" ""
Set-CellText $row7 3 "42" ""
Set-CellText $row7 4 "P8" ""
Set-CellText $row7 5 "P8" "O"

# --- Row 8: (was empty) -> Paper 2 / 72.0 / C4 ---------------------------
$row8 = $t.Rows.Item(8)
Set-CellText $row8 2 "" "Paper 2"
Set-CellText $row8 3 "" "72.0"
Set-CellText $row8 4 "" "C4"

# --- Row 9: (was empty) -> Paper 3 / 43.0 / P8 ---------------------------
$row9 = $t.Rows.Item(9)
Set-CellText $row9 2 "" "Paper 3"
Set-CellText $row9 3 "" "43.0"
Set-CellText $row9 4 "" "P8"

# --- Row 13: SUBMATH / 50 / P7 / P7 --------------------------------------
$row13 = $t.Rows.Item(13)
Set-CellText $row13 1 "SUBMATH" ""
Set-CellText $row13 3 "50" ""
Set-CellText $row13 4 "P7" ""
Set-CellText $row13 5 "P7" ""

# --- Row 14: GENERAL PAPER / 48 -> 33.0 / P8 -> F9 / P8 -> F9 ------------
$row14 = $t.Rows.Item(14)
Set-CellText $row14 3 "48" "33.0"
Set-CellText $row14 4 "P8" "F9"
Set-CellText $row14 5 "P8" "F9"
